$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "NIK" column header to "NPK"
$ws.Range("B1").Value = "NPK"
